# Apply cryptos list update (Thu Jan  4 22:34:24 UTC 2024) via GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "44.631.76"
$ws.Range('E2').Value = "  +4.04%  "
$ws.Range('D3').Value = "2.274.70"
$ws.Range('E3').Value = "  +2.70%  "
$ws.Range('E4').Value = "  +0.01%  "
$ws.Range('D5').Value = "'321.80"
$ws.Range('E5').Value = "  +1.73%  "
$ws.Range('D6').Value = "'106.10"
$ws.Range('E6').Value = "  +6.49%  "
$ws.Range('D7').Value = "'0.595"
$ws.Range('E7').Value = "  +0.83%  "
$ws.Range('E8').Value = "  +0.03%  "
$ws.Range('D9').Value = "'0.574"
$ws.Range('E9').Value = "  +2.49%  "
$ws.Range('D10').Value = "'38.68"
$ws.Range('E10').Value = "  +4.48%  "
$ws.Range('D11').Value = "'0.0844"
$ws.Range('E11').Value = "  +1.69%  "
$ws.Range('D12').Value = "'7.90"
$ws.Range('E12').Value = "  +2.29%  "
$ws.Range('E13').Value = "  +0.84%  "
$ws.Range('D14').Value = "'0.886"
$ws.Range('E14').Value = "  +3.23%  "
$ws.Range('D15').Value = "2.625.25"
$ws.Range('E15').Value = "  +2.67%  "
$ws.Range('D16').Value = "'14.64"
$ws.Range('E16').Value = "  +3.24%  "
$ws.Range('D17').Value = "2.286.86"
$ws.Range('E17').Value = "  +3.44%  "
$ws.Range('D18').Value = "44.542.66"
$ws.Range('E18').Value = "  +4.06%  "
$ws.Range('D19').Value = "'14.02"
$ws.Range('E19').Value = "  -4.39%  "
$ws.Range('D21').Value = "'6.55"
$ws.Range('E21').Value = "  +2.29%  "
$ws.Range('D22').Value = "'66.57"
$ws.Range('E22').Value = "  +2.05%  "
$ws.Range('D23').Value = "'3.22"
$ws.Range('E23').Value = "  +2.84%  "
$ws.Range('D24').Value = "'240.01"
$ws.Range('E24').Value = "  +1.87%  "
$ws.Range('D25').Value = "'2.22"
$ws.Range('E25').Value = "  +3.83%  "
$ws.Range('E26').Value = "  +0.04%  "
$ws.Range('D27').Value = "'10.20"
$ws.Range('E27').Value = "  +2.23%  "
$ws.Range('D28').Value = "'38.60"
$ws.Range('E28').Value = "  +12.24%  "
$ws.Range('E29').Value = "  +0.62%  "
$ws.Range('D30').Value = "'6.51"
$ws.Range('E30').Value = "  +3.25%  "
$ws.Range('D31').Value = "'20.69"
$ws.Range('E31').Value = "  +0.80%  "
$ws.Range('D32').Value = "'0.0888"
$ws.Range('E32').Value = "  -0.21%  "
$ws.Range('D33').Value = "'162.34"
$ws.Range('E33').Value = "  +4.33%  "
$ws.Range('D34').Value = "'2.79"
$ws.Range('E34').Value = "  +0.22%  "
$ws.Range('D36').Value = "'2.03"
$ws.Range('E36').Value = "  +5.10%  "
$ws.Range('D37').Value = "'3.16"
$ws.Range('E37').Value = "  +0.57%  "
$ws.Range('E38').Value = "  +0.35%  "
$ws.Range('D39').Value = "'3.93"
$ws.Range('E39').Value = "  +0.56%  "
$ws.Range('D40').Value = "'4.48"
$ws.Range('E40').Value = "  +1.36%  "
$ws.Range('D41').Value = "'15.62"
$ws.Range('E41').Value = "  +25.13%  "
$ws.Range('D42').Value = "'0.0329"
$ws.Range('E42').Value = "  +1.48%  "
$ws.Range('E43').Value = "  +0.20%  "
$ws.Range('D44').Value = "1.770.58"
$ws.Range('E44').Value = "  -7.25%  "
$ws.Range('B45').Value = "Algorand"
$ws.Range('C45').Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('D45').Value = "'0.210"
$ws.Range('E45').Value = "  +1.42%  "
$ws.Range('B46').Value = "BitcoinSV"
$ws.Range('C46').Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range('D46').Value = "'87.35"
$ws.Range('E46').Value = "  -1.13%  "
$ws.Range('D47').Value = "'5.48"
$ws.Range('E47').Value = "  +2.18%  "
$ws.Range('B48').Value = "MultiversX"
$ws.Range('C48').Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range('D48').Value = "'60.61"
$ws.Range('E48').Value = "  -0.33%  "
$ws.Range('B49').Value = "ordi"
$ws.Range('C49').Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range('D49').Value = "'76.00"
$ws.Range('E49').Value = "  -1.06%  "
$ws.Range('D50').Value = "'1.72"
$ws.Range('E50').Value = "  +7.85%  "
$ws.Range('D51').Value = "'104.59"
$ws.Range('E51').Value = "  +1.87%  "
